# Weekly driver report update for 2025-04-20
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) "Bad Drivers" table
#    Row 12 (Intel(R) Wi-Fi 6 AX201 160MHz - 22.250.10.1) drops out of the
#    report entirely this week, so remove that row first; this shifts the
#    Totals row (13 -> 12) and everything below it up by one.
# ---------------------------------------------------------------------------
$ws.Rows("12:12").Delete()

# Updated figures for the remaining bad-driver rows.
$ws.Range("C4").Value() = 3437
$ws.Range("D4").Value() = 86.8

$ws.Range("D5").Value() = 96

$ws.Range("C6").Value() = 872
$ws.Range("D6").Value() = 96.8

$ws.Range("C7").Value() = 2237
$ws.Range("D7").Value() = 98

# Row 8 now reports on driver 23.60.1.2 (was 22.170.2.1)
$ws.Range("A8").Value() = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.60.1.2"
$ws.Range("B8").Value() = 11
$ws.Range("C8").Value() = 476

$ws.Range("B9").Value() = 251
$ws.Range("C9").Value() = 8930

# Row 10 now reports on driver 22.170.2.1 (was 23.60.1.2)
$ws.Range("A10").Value() = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.170.2.1"
$ws.Range("B10").Value() = 509
$ws.Range("C10").Value() = 21814
$ws.Range("D10").Value() = 98.3

$ws.Range("B11").Value() = 300
$ws.Range("C11").Value() = 7018
$ws.Range("D11").Value() = 98.8

# Totals row (now row 12 after the delete above)
$ws.Range("B12").Value() = 1135
$ws.Range("C12").Value() = 45160

# ---------------------------------------------------------------------------
# 2) "Good Drivers" table
#    A brand-new top entry (23.100.0.4, refreshed sample count) is added,
#    pushing the header block down by one row and re-sorting the remaining
#    entries by driver vintage. Insert a row at 20 (just under the column
#    header, which is currently at row 19 once shifted) to make room.
# ---------------------------------------------------------------------------
$ws.Rows("20:20").Insert()

# Copy formatting from the row below (the old top data row, now at 21) so the
# new data row matches the rest of the table's styling.
$ws.Range("A21:E21").Copy()
$ws.Range("A20:E20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A20").Value() = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B20").Value() = 445055
$ws.Range("D20").Value() = 99.90000000000001
$ws.Range("E20").Formula() = "=""2024-11-10"""

# Re-populate rows 21-28 with the re-sorted data (driver vintage descending).
$ws.Range("A21").Value() = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1"
$ws.Range("B21").Value() = 10661
$ws.Range("D21").Value() = 100
$ws.Range("E21").Formula() = "=""2022-08-29"""

$ws.Range("A22").Value() = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3"
$ws.Range("B22").Value() = 14239
$ws.Range("D22").Value() = 100
$ws.Range("E22").Formula() = "=""2022-05-23"""

$ws.Range("A23").Value() = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1"
$ws.Range("B23").Value() = 265400
$ws.Range("D23").Value() = 99.90000000000001
$ws.Range("E23").Formula() = "=""2022-05-01"""

$ws.Range("A24").Value() = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B24").Value() = 77849
$ws.Range("D24").Value() = 99.90000000000001
$ws.Range("E24").Formula() = "=""2021-08-18"""

$ws.Range("A25").Value() = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B25").Value() = 34244
$ws.Range("D25").Value() = 100
$ws.Range("E25").Formula() = "=""2021-04-27"""

$ws.Range("A26").Value() = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B26").Value() = 59673
$ws.Range("D26").Value() = 100
$ws.Range("E26").Formula() = "=""2020-08-05"""

$ws.Range("A27").Value() = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B27").Value() = 113652
$ws.Range("D27").Value() = 100
$ws.Range("E27").Formula() = "=""2020-01-06"""

$ws.Range("A28").Value() = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B28").Value() = 56018
$ws.Range("D28").Value() = 100
$ws.Range("E28").Formula() = "=""2019-12-14"""

# Column C in this table is only ever a spacer column with no real content;
# clean up the placeholder values that the row insert/shift left behind so
# the data rows match the blank spacer cells used throughout the table.
$ws.Range("C21:C28").ClearContents()

# The row that used to be blank (29) picked up the old row 28's formatting
# when everything shifted down after the insert above; the table now ends
# at row 28, so fully clear row 29 back to blank again.
$ws.Range("A29:E29").Clear()

# ---------------------------------------------------------------------------
# 3) Sheet dimension shrinks by one row (A1:J34 -> A1:J33) since the bad
#    drivers table lost a row and the good drivers table gained one back.
#    Touch the bottom-right corner of the formatted (but otherwise unused)
#    region so the sheet's extent still reaches column J / row 33, matching
#    the original layout's reserved spacer area.
# ---------------------------------------------------------------------------
$ws.Range("J33").Font.Bold = $false
$ws.Range("A1").Select()
